# chore: adapt column header formatting to respective input file names
#
# The workbook's header row (row 1) contains the same 10 column names
# twice - once describing the "old"/"new" comparison side of an AHB diff.
# Rename those suffixes to the respective format-version they now refer
# to (_old -> _FV2410, _new -> _FV2504), wrap the sheet's data range in an
# Excel Table ("Table1") using those (renamed) headers as its column
# names, and freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns A1:U1 (21 columns), in order.
$headers = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410",
  "diff",
  "Segmentname_FV2504",
  "Segmentgruppe_FV2504",
  "Segment_FV2504",
  "Datenelement_FV2504",
  "Segment ID_FV2504",
  "Code_FV2504",
  "Qualifier_FV2504",
  "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504",
  "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table ("Table1") so the renamed
# headers become the table's column headers and an autofilter is shown.
$usedRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1) so it stays visible.
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
Write-Host "Applied header renames, Table1, and a frozen header row."
